$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "Илюха крут"
$ws.Range("B15").Value = "Илюха лох"

$ws.Rows.Item(16).Delete()

$ws.Range("B16").Select()
